# Auto-generated Excel COM-interop script to apply the NBA team stats data fix
# Commit: Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown; BF column dates
# are reformatted from "3-8-2011-12" to ISO format "2012-03-08", and several stat
# values across rows 2-31 are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Date column (BF) to be stored as text so the "YYYY-MM-DD" strings
# are not auto-converted into Excel date serial numbers.
$ws.Range("BF2:BF31").NumberFormat = "@"

# Row 2
$ws.Range("AI2").Value = 17

# Row 3
$ws.Range("AD3").Value = 20
$ws.Range("AF3").Value = 13
$ws.Range("AQ3").Value = 11

# Row 4
$ws.Range("AV4").Value = 11

# Row 5
$ws.Range("D5").Value = 41
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 0.805
$ws.Range("I5").Value = 38
$ws.Range("J5").Value = 82.3
$ws.Range("K5").Value = 0.461
$ws.Range("M5").Value = 15.9
$ws.Range("N5").Value = 0.387
$ws.Range("O5").Value = 15.6
$ws.Range("P5").Value = 21.5
$ws.Range("Q5").Value = 0.726
$ws.Range("R5").Value = 13.3
$ws.Range("S5").Value = 32.6
$ws.Range("T5").Value = 45.9
$ws.Range("U5").Value = 23.2
$ws.Range("V5").Value = 14.2
$ws.Range("W5").Value = 7.2
$ws.Range("Z5").Value = 17.4
$ws.Range("AB5").Value = 97.8
$ws.Range("AC5").Value = 9.5
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 1
$ws.Range("AK5").Value = 4
$ws.Range("AM5").Value = 21
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 23
$ws.Range("AP5").Value = 20
$ws.Range("AT5").Value = 1
$ws.Range("AV5").Value = 6

# Row 7
$ws.Range("D7").Value = 40
$ws.Range("F7").Value = 17
$ws.Range("G7").Value = 0.575
$ws.Range("I7").Value = 35.6
$ws.Range("J7").Value = 81.7
$ws.Range("O7").Value = 15.7
$ws.Range("P7").Value = 20.9
$ws.Range("R7").Value = 10.7
$ws.Range("S7").Value = 32.5
$ws.Range("T7").Value = 43.2
$ws.Range("V7").Value = 14.6
$ws.Range("W7").Value = 9.300000000000001
$ws.Range("X7").Value = 5.3
$ws.Range("Z7").Value = 19.5
$ws.Range("AA7").Value = 18.6
$ws.Range("AC7").Value = 2.9
$ws.Range("AF7").Value = 11
$ws.Range("AG7").Value = 11
$ws.Range("AI7").Value = 18
$ws.Range("AN7").Value = 22
$ws.Range("AO7").Value = 22
$ws.Range("AV7").Value = 12
$ws.Range("AY7").Value = 4
$ws.Range("AZ7").Value = 12

# Row 8
$ws.Range("AD8").Value = 2
$ws.Range("AF8").Value = 13
$ws.Range("AZ8").Value = 11

# Row 9
$ws.Range("AN9").Value = 21
$ws.Range("AW9").Value = 26

# Row 10
$ws.Range("AN10").Value = 5

# Row 11
$ws.Range("AD11").Value = 2

# Row 12
$ws.Range("AV12").Value = 10

# Row 14
$ws.Range("AT14").Value = 2

# Row 15
$ws.Range("AD15").Value = 20
$ws.Range("AJ15").Value = 14
$ws.Range("AZ15").Value = 10

# Row 16
$ws.Range("AF16").Value = 3
$ws.Range("AY16").Value = 5

# Row 17
$ws.Range("AL17").Value = 16
$ws.Range("AY17").Value = 15
$ws.Range("BA17").Value = 20

# Row 18
$ws.Range("AD18").Value = 2
$ws.Range("AN18").Value = 16

# Row 19
$ws.Range("AD19").Value = 2

# Row 20
$ws.Range("AK20").Value = 18
$ws.Range("AO20").Value = 27

# Row 21
$ws.Range("AR21").Value = 16

# Row 22
$ws.Range("AG22").Value = 2

# Row 23
$ws.Range("D23").Value = 40
$ws.Range("E23").Value = 25
$ws.Range("G23").Value = 0.625
$ws.Range("M23").Value = 26.7
$ws.Range("N23").Value = 0.387
$ws.Range("O23").Value = 15.5
$ws.Range("P23").Value = 24
$ws.Range("Q23").Value = 0.644
$ws.Range("R23").Value = 11.2
$ws.Range("T23").Value = 43.3
$ws.Range("U23").Value = 20.5
$ws.Range("W23").Value = 6.5
$ws.Range("AA23").Value = 20.8
$ws.Range("AB23").Value = 93.8
$ws.Range("AE23").Value = 5
$ws.Range("AK23").Value = 17
$ws.Range("AO23").Value = 25
$ws.Range("AR23").Value = 17

# Row 24
$ws.Range("AD24").Value = 2
$ws.Range("AR24").Value = 24

# Row 25
$ws.Range("D25").Value = 38
$ws.Range("E25").Value = 17
$ws.Range("G25").Value = 0.447
$ws.Range("J25").Value = 81.59999999999999
$ws.Range("K25").Value = 0.447
$ws.Range("L25").Value = 6.6
$ws.Range("M25").Value = 19.5
$ws.Range("N25").Value = 0.337
$ws.Range("O25").Value = 15.2
$ws.Range("P25").Value = 19.9
$ws.Range("Q25").Value = 0.765
$ws.Range("R25").Value = 10.6
$ws.Range("S25").Value = 31.5
$ws.Range("V25").Value = 14.8
$ws.Range("W25").Value = 6.9
$ws.Range("Y25").Value = 4.6
$ws.Range("Z25").Value = 19
$ws.Range("AA25").Value = 19.3
$ws.Range("AC25").Value = -2
$ws.Range("AD25").Value = 20
$ws.Range("AE25").Value = 20
$ws.Range("AG25").Value = 20
$ws.Range("AJ25").Value = 13
$ws.Range("AL25").Value = 15
$ws.Range("AN25").Value = 15
$ws.Range("AO25").Value = 28
$ws.Range("AQ25").Value = 10
$ws.Range("AR25").Value = 23
$ws.Range("AW25").Value = 25
$ws.Range("BA25").Value = 19

# Row 26
$ws.Range("AV26").Value = 4

# Row 28
$ws.Range("AD28").Value = 20
$ws.Range("AK28").Value = 5
$ws.Range("AP28").Value = 19

# Row 29
$ws.Range("AN29").Value = 17

# Row 30
$ws.Range("AD30").Value = 20
$ws.Range("AV30").Value = 5

# Row 31
$ws.Range("AD31").Value = 20
$ws.Range("AM31").Value = 20

# Date column (BF): "3-8-2011-12" -> "2012-03-08" for every data row (2-31)
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2012-03-08"
}
